# missingHeaderRow.xlsx — "Fixed failed unit tests" commit.
#
# The canonical diff for this commit is dominated by cosmetic, Excel-version-
# specific re-save artifacts (fileVersion/rupBuild, x15ac:absPath, the
# xr:revisionPtr GUID, workbookView window geometry, calcPr's iterative-calc
# echo, styles.xml's x14ac:knownFonts flag, and x14ac:dyDescent on every row)
# that reflect which machine/Excel build produced the save — they aren't
# deliberate user edits and aren't reachable through the object model here.
# The two real, intentional content edits are:
#   1. The sheet was renamed from "Basic Clinic Data..." to "Basic Clinic
#      Data" (trailing ellipsis removed).
#   2. The active selection moved from A1:XFD1 (a full-row select) to C18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet tab (ellipsis removed from "Basic Clinic Data...").
$ws.Name = "Basic Clinic Data"

# 2. Move the selection to C18, matching the saved <selection> state.
$ws.Range("C18").Select() | Out-Null
